$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2625   # H40
$ws.Cells.Item(40, 9).Value = 3028.5715   # I40
$ws.Cells.Item(40, 10).Value = 2060   # J40
$ws.Cells.Item(40, 11).Value = 3028.5715   # K40
$ws.Cells.Item(40, 12).Value = 2060   # L40
$ws.Cells.Item(40, 13).Value = -2853.5715   # M40
$ws.Cells.Item(40, 14).Value = -2410   # N40
$ws.Cells.Item(113, 8).Value = 2774.5386   # H113
$ws.Cells.Item(113, 9).Value = 2696.75   # I113
$ws.Cells.Item(113, 11).Value = 2696.75   # K113
$ws.Cells.Item(113, 13).Value = 557.25   # M113
$ws.Cells.Item(134, 8).Value = 60000   # H134
$ws.Cells.Item(134, 10).Value = 60000   # J134
$ws.Cells.Item(134, 12).Value = 60000   # L134
$ws.Cells.Item(134, 14).Value = -70140   # N134
$ws.Cells.Item(137, 8).Value = 1049.2941   # H137
$ws.Cells.Item(137, 9).Value = 768.7241   # I137
$ws.Cells.Item(137, 10).Value = 1652   # J137
$ws.Cells.Item(137, 11).Value = 2306.1723   # K137
$ws.Cells.Item(137, 12).Value = 4956   # L137
$ws.Cells.Item(137, 13).Value = 243.8276999999998   # M137
$ws.Cells.Item(137, 14).Value = -10056   # N137
$ws.Cells.Item(138, 8).Value = 588764.4   # H138
$ws.Cells.Item(138, 9).Value = 688.9697   # I138
$ws.Cells.Item(138, 10).Value = 1062093.4   # J138
$ws.Cells.Item(138, 11).Value = 2066.9091   # K138
$ws.Cells.Item(138, 12).Value = 3186280.2   # L138
$ws.Cells.Item(138, 13).Value = 3073.0909   # M138
$ws.Cells.Item(138, 14).Value = -3196560.2   # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3943.2134   # H32
$ws.Cells.Item(32, 9).Value = 3548.791   # I32
$ws.Cells.Item(32, 10).Value = 7246.5   # J32
$ws.Cells.Item(32, 11).Value = 3548.791   # K32
$ws.Cells.Item(32, 12).Value = 7246.5   # L32
$ws.Cells.Item(32, 13).Value = -3261.791   # M32
$ws.Cells.Item(32, 14).Value = -7820.5   # N32
$ws.Cells.Item(61, 8).Value = 14706829   # H61
$ws.Cells.Item(61, 9).Value = 16949964   # I61
$ws.Cells.Item(61, 10).Value = 1834.8889   # J61
$ws.Cells.Item(61, 11).Value = 16949964   # K61
$ws.Cells.Item(61, 12).Value = 1834.8889   # L61
$ws.Cells.Item(61, 13).Value = -16949752   # M61
$ws.Cells.Item(61, 14).Value = -2258.8889   # N61
$ws.Cells.Item(122, 8).Value = 3037.3333   # H122
$ws.Cells.Item(122, 9).Value = 3037.3333   # I122
$ws.Cells.Item(122, 10).Value = 0   # J122
$ws.Cells.Item(122, 11).Value = 9111.999899999999   # K122
$ws.Cells.Item(122, 12).Value = 0   # L122
$ws.Cells.Item(122, 13).Value = -6661.999899999999   # M122
$ws.Cells.Item(122, 14).ClearContents()   # N122
$ws.Cells.Item(132, 8).Value = 1344.0702   # H132
$ws.Cells.Item(132, 9).Value = 1257.84   # I132
$ws.Cells.Item(132, 10).Value = 1960   # J132
$ws.Cells.Item(132, 11).Value = 3773.52   # K132
$ws.Cells.Item(132, 12).Value = 5880   # L132
$ws.Cells.Item(132, 13).Value = -1243.52   # M132
$ws.Cells.Item(132, 14).Value = -10940   # N132
$ws.Cells.Item(136, 8).Value = 14706829   # H136
$ws.Cells.Item(136, 9).Value = 16949964   # I136
$ws.Cells.Item(136, 10).Value = 1834.8889   # J136
$ws.Cells.Item(136, 11).Value = 50849892   # K136
$ws.Cells.Item(136, 12).Value = 5504.6667   # L136
$ws.Cells.Item(136, 13).Value = -50847342   # M136
$ws.Cells.Item(136, 14).Value = -10604.6667   # N136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2809.8635   # H134
$ws.Cells.Item(134, 9).Value = 775.0172   # I134
$ws.Cells.Item(134, 10).Value = 17562.5   # J134
$ws.Cells.Item(134, 11).Value = 2325.0516   # K134
$ws.Cells.Item(134, 12).Value = 52687.5   # L134
$ws.Cells.Item(134, 13).Value = 209.9484000000002   # M134
$ws.Cells.Item(134, 14).Value = -57757.5   # N134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 47780   # H20
$ws.Cells.Item(20, 10).Value = 47780   # J20
$ws.Cells.Item(20, 12).Value = 47780   # L20
$ws.Cells.Item(20, 14).Value = -48252   # N20
$ws.Cells.Item(30, 8).Value = 47780   # H30
$ws.Cells.Item(30, 10).Value = 47780   # J30
$ws.Cells.Item(30, 12).Value = 47780   # L30
$ws.Cells.Item(30, 14).Value = -47962   # N30
$ws.Cells.Item(31, 8).Value = 1975.5   # H31
$ws.Cells.Item(31, 9).Value = 2090.9092   # I31
$ws.Cells.Item(31, 10).Value = 1552.3334   # J31
$ws.Cells.Item(31, 11).Value = 2090.9092   # K31
$ws.Cells.Item(31, 12).Value = 1552.3334   # L31
$ws.Cells.Item(31, 13).Value = -1795.9092   # M31
$ws.Cells.Item(31, 14).Value = -2142.3334   # N31
$ws.Cells.Item(34, 8).Value = 1975.5   # H34
$ws.Cells.Item(34, 9).Value = 2090.9092   # I34
$ws.Cells.Item(34, 10).Value = 1552.3334   # J34
$ws.Cells.Item(34, 11).Value = 2090.9092   # K34
$ws.Cells.Item(34, 12).Value = 1552.3334   # L34
$ws.Cells.Item(34, 13).Value = -1888.9092   # M34
$ws.Cells.Item(34, 14).Value = -1956.3334   # N34
$ws.Cells.Item(99, 8).Value = 1542.6666   # H99
$ws.Cells.Item(99, 9).Value = 1599.1666   # I99
$ws.Cells.Item(99, 10).Value = 1316.6666   # J99
$ws.Cells.Item(99, 11).Value = 1599.1666   # K99
$ws.Cells.Item(99, 12).Value = 1316.6666   # L99
$ws.Cells.Item(99, 13).Value = -101.1666   # M99
$ws.Cells.Item(99, 14).Value = -4312.6666   # N99
$ws.Cells.Item(122, 8).Value = 1221.2632   # H122
$ws.Cells.Item(122, 9).Value = 1107.5   # I122
$ws.Cells.Item(122, 10).Value = 1539.8   # J122
$ws.Cells.Item(122, 11).Value = 3322.5   # K122
$ws.Cells.Item(122, 12).Value = 4619.4   # L122
$ws.Cells.Item(122, 13).Value = -872.5   # M122
$ws.Cells.Item(122, 14).Value = -9519.4   # N122
$ws.Cells.Item(126, 8).Value = 1542.6666   # H126
$ws.Cells.Item(126, 9).Value = 1599.1666   # I126
$ws.Cells.Item(126, 10).Value = 1316.6666   # J126
$ws.Cells.Item(126, 11).Value = 4797.4998   # K126
$ws.Cells.Item(126, 12).Value = 3949.9998   # L126
$ws.Cells.Item(126, 13).Value = -2327.4998   # M126
$ws.Cells.Item(126, 14).Value = -8889.9998   # N126
$ws.Cells.Item(128, 8).Value = 47780   # H128
$ws.Cells.Item(128, 10).Value = 47780   # J128
$ws.Cells.Item(128, 12).Value = 47780   # L128
$ws.Cells.Item(128, 14).Value = -57740   # N128
$ws.Cells.Item(129, 8).Value = 35555.715   # H129
$ws.Cells.Item(129, 9).Value = 10000   # I129
$ws.Cells.Item(129, 10).Value = 39815   # J129
$ws.Cells.Item(129, 11).Value = 10000   # K129
$ws.Cells.Item(129, 12).Value = 39815   # L129
$ws.Cells.Item(129, 13).Value = -5000   # M129
$ws.Cells.Item(129, 14).Value = -49815   # N129
$ws.Cells.Item(132, 8).Value = 1432.6031   # H132
$ws.Cells.Item(132, 9).Value = 1324.5927   # I132
$ws.Cells.Item(132, 10).Value = 2080.6667   # J132
$ws.Cells.Item(132, 11).Value = 3973.7781   # K132
$ws.Cells.Item(132, 12).Value = 6242.000100000001   # L132
$ws.Cells.Item(132, 13).Value = -1443.7781   # M132
$ws.Cells.Item(132, 14).Value = -11302.0001   # N132

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 98.833336   # H2
$ws.Cells.Item(2, 10).Value = 133.16667   # J2
$ws.Cells.Item(2, 12).Value = 799.0000200000001   # L2
$ws.Cells.Item(2, 14).Value = -1025.00002   # N2
$ws.Cells.Item(125, 8).Value = 6152.5   # H125
$ws.Cells.Item(125, 10).Value = 6169.4443   # J125
$ws.Cells.Item(125, 12).Value = 18508.3329   # L125
$ws.Cells.Item(125, 14).Value = -28348.3329   # N125
$ws.Cells.Item(131, 8).Value = 19232054   # H131
$ws.Cells.Item(131, 9).Value = 166667150   # I131
$ws.Cells.Item(131, 10).Value = 1389.2826   # J131
$ws.Cells.Item(131, 11).Value = 500001450   # K131
$ws.Cells.Item(131, 12).Value = 4167.8478   # L131
$ws.Cells.Item(131, 13).Value = -499996410   # M131
$ws.Cells.Item(131, 14).Value = -14247.8478   # N131
$ws.Cells.Item(137, 8).Value = 19741978   # H137
$ws.Cells.Item(137, 9).Value = 41668310   # I137
$ws.Cells.Item(137, 10).Value = 8273.85   # J137
$ws.Cells.Item(137, 11).Value = 125004930   # K137
$ws.Cells.Item(137, 12).Value = 24821.55   # L137
$ws.Cells.Item(137, 13).Value = -124999830   # M137
$ws.Cells.Item(137, 14).Value = -35021.55   # N137
$ws.Cells.Item(139, 8).Value = 1654.3611   # H139
$ws.Cells.Item(139, 9).Value = 1647.4814   # I139
$ws.Cells.Item(139, 10).Value = 1675   # J139
$ws.Cells.Item(139, 11).Value = 4942.4442   # K139
$ws.Cells.Item(139, 12).Value = 5025   # L139
$ws.Cells.Item(139, 13).Value = 197.5558000000001   # M139
$ws.Cells.Item(139, 14).Value = -15305   # N139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2057.0476   # H122
$ws.Cells.Item(122, 9).Value = 1512.4375   # I122
$ws.Cells.Item(122, 10).Value = 3799.8   # J122
$ws.Cells.Item(122, 11).Value = 4537.3125   # K122
$ws.Cells.Item(122, 12).Value = 11399.4   # L122
$ws.Cells.Item(122, 13).Value = -2087.3125   # M122
$ws.Cells.Item(122, 14).Value = -16299.4   # N122
$ws.Cells.Item(126, 8).Value = 2331.8333   # H126
$ws.Cells.Item(126, 10).Value = 2669   # J126
$ws.Cells.Item(126, 12).Value = 8007   # L126
$ws.Cells.Item(126, 14).Value = -12947   # N126
$ws.Cells.Item(132, 8).Value = 1339.3611   # H132
$ws.Cells.Item(132, 9).Value = 924.23334   # I132
$ws.Cells.Item(132, 11).Value = 2772.70002   # K132
$ws.Cells.Item(132, 13).Value = -242.7000200000002   # M132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2390.64   # H40
$ws.Cells.Item(40, 9).Value = 1641.7826   # I40
$ws.Cells.Item(40, 11).Value = 1641.7826   # K40
$ws.Cells.Item(40, 13).Value = -1505.7826   # M40
$ws.Cells.Item(46, 8).Value = 1682.8572   # H46
$ws.Cells.Item(46, 9).Value = 1390   # I46
$ws.Cells.Item(46, 10).Value = 1800   # J46
$ws.Cells.Item(46, 11).Value = 1390   # K46
$ws.Cells.Item(46, 12).Value = 1800   # L46
$ws.Cells.Item(46, 13).Value = -1202   # M46
$ws.Cells.Item(46, 14).Value = -2176   # N46
$ws.Cells.Item(132, 8).Value = 19364.965   # H132
$ws.Cells.Item(132, 9).Value = 1193.1   # I132
$ws.Cells.Item(132, 10).Value = 64794.625   # J132
$ws.Cells.Item(132, 11).Value = 3579.3   # K132
$ws.Cells.Item(132, 12).Value = 194383.875   # L132
$ws.Cells.Item(132, 13).Value = -1049.3   # M132
$ws.Cells.Item(132, 14).Value = -199443.875   # N132
$ws.Cells.Item(136, 8).Value = 1187.1177   # H136
$ws.Cells.Item(136, 9).Value = 1068.7333   # I136
$ws.Cells.Item(136, 11).Value = 3206.199900000001   # K136
$ws.Cells.Item(136, 13).Value = -656.1999000000005   # M136

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 495.42856   # H107
$ws.Cells.Item(107, 9).Value = 459.63635   # I107
$ws.Cells.Item(107, 11).Value = 1378.90905   # K107
$ws.Cells.Item(107, 13).Value = 541.09095   # M107
$ws.Cells.Item(125, 8).Value = 40292.145   # H125
$ws.Cells.Item(125, 10).Value = 40292.145   # J125
$ws.Cells.Item(125, 12).Value = 40292.145   # L125
$ws.Cells.Item(125, 14).Value = -50132.145   # N125
$ws.Cells.Item(132, 8).Value = 2132.0466   # H132
$ws.Cells.Item(132, 9).Value = 2337.2222   # I132
$ws.Cells.Item(132, 10).Value = 1076.8572   # J132
$ws.Cells.Item(132, 11).Value = 7011.6666   # K132
$ws.Cells.Item(132, 12).Value = 3230.5716   # L132
$ws.Cells.Item(132, 13).Value = -4481.6666   # M132
$ws.Cells.Item(132, 14).Value = -8290.571599999999   # N132
$ws.Cells.Item(136, 8).Value = 619.51514   # H136
$ws.Cells.Item(136, 9).Value = 425.9091   # I136
$ws.Cells.Item(136, 10).Value = 1006.7273   # J136
$ws.Cells.Item(136, 11).Value = 1277.7273   # K136
$ws.Cells.Item(136, 12).Value = 3020.1819   # L136
$ws.Cells.Item(136, 13).Value = 1272.2727   # M136
$ws.Cells.Item(136, 14).Value = -8120.1819   # N136
